# Inserts a new data row at row 357 (pushing the existing rows 357-385 down
# to 358-386) and populates the new row with its own data, matching the
# "Hortaliza, Macroferia Regional de Talca - Repollo" weekly update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 357:385 down to 358:386, leaving a blank row 357 behind.
$ws.Rows(357).Insert()

# Populate the newly inserted row 357 with the new record.
$ws.Range("A357").Value = 5
$ws.Range("B357").Value = "Macroferia Regional de Talca"
$ws.Range("C357").Value = "Maule"
$ws.Range("D357").Value = 44826
$ws.Range("E357").Value = 7
$ws.Range("F357").Value = 100112006
$ws.Range("G357").Value = "Repollo"
$ws.Range("H357").Value = "Crespo record"
$ws.Range("I357").Value = "Primera"
$ws.Range("J357").Value = 3000
$ws.Range("K357").Value = 2000
$ws.Range("L357").Value = 2000
$ws.Range("M357").Value = 2000
$ws.Range("N357").Value = '$/unidad'
$ws.Range("O357").Value = "Provincia del Elquí"
$ws.Range("P357").Value = 2000
$ws.Range("Q357").Value = 1
$ws.Range("R357").Value = "Hortaliza"
